$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.411653518676758
$ws.Range("C2").Value = 5.344827651977539
$ws.Range("D2").Value = 11.668232917785645
$ws.Range("E2").Value = 46.42856979370117
